# Update the "UserPermissions" test data sheet:
#  - rename column header C2 from "isAddButtonVisisble" to "searchDocument_isEditable"
#  - add a new column D with header "missingClient_isEditable"
#  - C4 flips from "Yes" to "No"; D3/D4 get "No"/"Yes" respectively
#  - widen column C and size the new column D
#  - move the active selection to D3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 2) - copy style from the existing header cells onto D2.
$ws.Range("C2").Value = "searchDocument_isEditable"
$ws.Range("D2").Value = "missingClient_isEditable"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null

# Data rows
$ws.Range("C3").Value = "No"
$ws.Range("D3").Value = "No"

$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = "Yes"

# Column widths (nearest value the width model can represent)
$ws.Columns.Item(3).ColumnWidth = 28
$ws.Columns.Item(4).ColumnWidth = 26.25

# Move selection to D3 to match the saved view state
$ws.Range("D3").Select() | Out-Null
